$wb = $excel.ActiveWorkbook

# Update the "想去人数" (want-to-go count) column on both the "展览"
# sheet and the "全部类型" sheet, which mirror the same rows.
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 117
    $ws.Range("F3").Value = 313
}
